$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header change
$ws.Range("B1").Value = "accession_number"

# Row 2
$ws.Range("A2").Value = "PRAX"
$ws.Range("E2").Value = "Phase 2a"
$ws.Range("F2").Value = "NCT05021978"
$ws.Range("G2").Value = "not specified"
$ws.Range("H2").Value = "A"
$ws.Range("I2").Value = "'2022"
$ws.Range("J2").Value = "Demonstrated positive results in the study."
$ws.Range("K2").Value = "Positive"
$ws.Range("O2").Value = "not specified"
$ws.Range("Q2").Value = "not specified"
$ws.Range("V2").Value = "positive results"
$ws.Range("W2").Value = "not specified"
$ws.Range("X2").Value = "Completed"

# Row 3
$ws.Range("A3").Value = "PRAX"
$ws.Range("E3").Value = "Phase 2b"
$ws.Range("F3").Value = "Essential 1"
$ws.Range("G3").Value = "Approximately 130"
$ws.Range("H3").Value = "A"
$ws.Range("I3").Value = "2023Q1"
$ws.Range("J3").Value = "Topline results from the Essential 1 study were announced. The study is a multi-center, randomized, double-blind, placebo-controlled, dose-range finding clinical trial evaluating the efficacy, safety and tolerability of once-daily treatment of ulixacaltamide compared to placebo."
$ws.Range("O3").Value = "multi-center, randomized, double-blind, placebo-controlled, dose-range finding"
$ws.Range("Q3").Value = "Placebo"
$ws.Range("S3").Value = "not specified"
$ws.Range("X3").Value = "Completed"

# Row 4
$ws.Range("A4").Value = "PRAX"
$ws.Range("D4").Value = "Parkinson’s disease"
$ws.Range("J4").Value = "A randomized, double-blind, placebo-controlled proof of concept trial planned to evaluate the efficacy, safety, and tolerability of ulixacaltamide as a non-dopaminergic treatment for motor symptoms of PD. Primary endpoint is change in UPDRS Part III motor examination score in the OFF state."
$ws.Range("V4").Value = "non-dopaminergic treatment for motor symptoms"
$ws.Range("X4").Value = "not specified"

# Row 5
$ws.Range("A5").Value = "PRAX"
$ws.Range("D5").Value = "Essential Tremor"
$ws.Range("E5").Value = "Phase 3"
$ws.Range("F5").Value = "Essential 3 (Study 1 interim analysis)"
$ws.Range("G5").Value = "N=400"
$ws.Range("I5").Value = "2025Q1"
$ws.Range("J5").Value = "Results of a pre-planned interim analysis of Study 1 of the Essential 3 clinical program were shared in February 2025. The IDMC recommended stopping for futility due to results being unlikely to meet the primary efficacy endpoint. The company decided to continue both studies to completion."
$ws.Range("K5").Value = "Futility"
$ws.Range("O5").Value = "decentralized, multi-study, 12-week parallel design, placebo-controlled study"
$ws.Range("Q5").Value = "Placebo"
$ws.Range("W5").Value = "Interim"
$ws.Range("X5").Value = "Ongoing"

# Row 6
$ws.Range("A6").Value = "PRAX"
$ws.Range("D6").Value = "Essential Tremor"
$ws.Range("F6").Value = "Essential 3 (Study 1 and Study 2 combined)"
$ws.Range("G6").Value = "'600"
$ws.Range("I6").Value = "2025Q3"
$ws.Range("J6").Value = "Topline results for both Study 1 and Study 2 of the Essential 3 program are expected. A decision about whether the data supports NDA submission will be made after analyzing the final results."
$ws.Range("M6").Value = "Potential NDA submission in 2025"
$ws.Range("O6").Value = "decentralized, multi-study, clinical trial evaluating safety and efficacy of 60 mg of ulixacaltamide in ET. Includes a 12-week parallel design, placebo-controlled study (Study 1) and a 12-week randomized withdrawal study (Study 2), with a long-term safety study (LTSS). Uses mADL11 as the primary endpoint."
$ws.Range("Q6").Value = "Placebo"
$ws.Range("S6").Value = "NDA"
$ws.Range("X6").Value = "Ongoing"

# Delete rows 7, 8, 9 (the trailing rows removed entirely)
$ws.Range("A7:X9").EntireRow.Delete()
